# Rating1.xlsx — "Added original huang experiment folder"
#
# The image-instructions folder referenced by the rating sheet moved from
# "Instructions_EN/" to a local "Rating\" folder (the author's original
# huang-experiment layout), so the two picture-path cells are updated to
# match. The active selection/view state is nudged as it was left after the
# edit, and the row heights / column width are touched up to mirror the
# resaved layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell A2 / A3 hold the relative image paths used by the rating instructions.
$ws.Range("A2").Value = "Rating\ratingCS+1.png"
$ws.Range("A3").Value = "Rating\ratingCS-1.png"

# Row heights for the header + two data rows were nudged slightly.
$ws.Rows.Item(1).RowHeight = 15.6
$ws.Rows.Item(2).RowHeight = 15.6
$ws.Rows.Item(3).RowHeight = 15.6

# Column B width was nudged slightly narrower.
$ws.Columns.Item(2).ColumnWidth = 16.3

# Selection left on A6 after the edit.
$ws.Range("A6").Select()
